{"js": "const pairs = [\n  [\"2023-12-13 Wednesday\", \"2023-12-14 Thursday\"],\n  [\"50+9=59\", \"29+11=40\"],\n  [\"51-15=36\", \"33-2=31\"],\n  [\"52-23=29\", \"24+3=27\"],\n  [\"68-46=22\", \"60-26=34\"],\n  [\"27-12=15\", \"63-25=38\"],\n  [\"92-35=57\", \"40+42=82\"],\n  [\"31-25=6\", \"32-12=20\"],\n  [\"76+22=98\", \"3+49=52\"],\n  [\"88-53=35\", \"14-5=9\"],\n  [\"20+1=21\", \"40+43=83\"],\n  [\"27+61=88\", \"50+6=56\"],\n  [\"30+6=36\", \"62-28=34\"],\n  [\"49+6=55\", \"20-2=18\"],\n  [\"87-0=87\", \"82-5=77\"],\n  [\"80-29=51\", \"51-33=18\"],\n  [\"67+19=86\", \"42-38=4\"],\n  [\"91-51=40\", \"15+27=42\"],\n  [\"50-4=46\", \"5+59=64\"],\n  [\"84-0=84\", \"13+66=79\"],\n  [\"22+17=39\", \"36+31=67\"],\n  [\"88-59=29\", \"67+30=97\"],\n  [\"18-17=1\", \"23+62=85\"],\n  [\"73+26=99\", \"86-37=49\"],\n  [\"18+50=68\", \"58+38=96\"],\n  [\"71-70=1\", \"34-7=27\"],\n  [\"36+44=80\", \"37+34=71\"],\n  [\"25-4=21\", \"38+25=63\"],\n  [\"12+76=88\", \"24+60=84\"],\n  [\"60-14=46\", \"83-70=13\"],\n  [\"65+5=70\", \"4+57=61\"],\n  [\"36+4=40\", \"78-0=78\"],\n  [\"24+39=63\", \"26+38=64\"],\n  [\"46-45=1\", \"90-22=68\"],\n  [\"73-47=26\", \"89+1=90\"],\n  [\"61-59=2\", \"88-70=18\"],\n  [\"45-6=39\", \"3+66=69\"],\n  [\"12-10=2\", \"85-72=13\"],\n  [\"0+64=64\", \"9+62=71\"],\n  [\"60-38=22\", \"77+18=95\"],\n  [\"72-51=21\", \"7+71=78\"],\n  [\"9+81=90\", \"10-9=1\"],\n  [\"52+3=55\", \"72-29=43\"],\n  [\"20-6=14\", \"22+44=66\"],\n  [\"73-30=43\", \"53-29=24\"],\n  [\"32-24=8\", \"46+36=82\"],\n  [\"32-25=7\", \"24+18=42\"],\n  [\"14+19=33\", \"41+27=68\"],\n  [\"76-64=12\", \"90-15=75\"],\n  [\"73+2=75\", \"20+71=91\"],\n  [\"42-32=10\", \"25+34=59\"],\n  [\"62+21=83\", \"3+81=84\"],\n  [\"87-56=31\", \"49-49=0\"],\n  [\"23-6=17\", \"74-19=55\"],\n  [\"38+24=62\", \"51+28=79\"],\n  [\"98-98=0\", \"34-16=18\"],\n  [\"58-8=50\", \"99-23=76\"],\n  [\"45-8=37\", \"98-76=22\"],\n  [\"52-22=30\", \"90-80=10\"],\n  [\"98-68=30\", \"40-30=10\"],\n  [\"98-8=90\", \"97-33=64\"],\n  [\"50+14=64\", \"13+52=65\"],\n  [\"2+30=32\", \"68-10=58\"],\n  [\"7+58=65\", \"75-13=62\"],\n  [\"43-31=12\", \"40-23=17\"],\n  [\"68-32=36\", \"22+18=40\"],\n  [\"99-93=6\", \"51+24=75\"],\n  [\"49-31=18\", \"59-4=55\"],\n  [\"48+14=62\", \"33+24=57\"],\n  [\"54-1=53\", \"18+67=85\"],\n  [\"24+16=40\", \"72+25=97\"],\n  [\"73-17=56\", \"21-16=5\"],\n  [\"72-56=16\", \"71-4=67\"],\n  [\"62-54=8\", \"46+47=93\"],\n  [\"31-14=17\", \"12+65=77\"],\n  [\"27+36=63\", \"95-42=53\"],\n  [\"5+34=39\", \"30+15=45\"],\n  [\"25+22=47\", \"12+16=28\"],\n  [\"84-30=54\", \"97-95=2\"],\n  [\"51-19=32\", \"85-27=58\"],\n  [\"20-19=1\", \"83+14=97\"],\n  [\"32+11=43\", \"14+42=56\"],\n  [\"78+10=88\", \"3+12=15\"],\n  [\"4+83=87\", \"83-56=27\"],\n  [\"5+62=67\", \"96-83=13\"],\n  [\"47+36=83\", \"68+21=89\"],\n  [\"13+53=66\", \"23+29=52\"],\n  [\"20+77=97\", \"51-12=39\"],\n  [\"67-21=46\", \"39-8=31\"],\n  [\"4+40=44\", \"95-45=50\"],\n  [\"81-47=34\", \"43+8=51\"],\n  [\"72-55=17\", \"8+71=79\"],\n  [\"82-16=66\", \"29+51=80\"],\n  [\"47-14=33\", \"10+24=34\"],\n  [\"78-65=13\", \"94-6=88\"],\n  [\"84-40=44\", \"99-37=62\"],\n  [\"63-34=29\", \"13+15=28\"],\n  [\"4+20=24\", \"79+12=91\"],\n  [\"45-42=3\", \"80-11=69\"],\n  [\"77+16=93\", \"21+60=81\"],\n  [\"99-24=75\", \"25+14=39\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWildcards: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2023-12-13 Wednesday', '2023-12-14 Thursday')\n    ,@('50+9=59', '29+11=40')\n    ,@('51-15=36', '33-2=31')\n    ,@('52-23=29', '24+3=27')\n    ,@('68-46=22', '60-26=34')\n    ,@('27-12=15', '63-25=38')\n    ,@('92-35=57', '40+42=82')\n    ,@('31-25=6', '32-12=20')\n    ,@('76+22=98', '3+49=52')\n    ,@('88-53=35', '14-5=9')\n    ,@('20+1=21', '40+43=83')\n    ,@('27+61=88', '50+6=56')\n    ,@('30+6=36', '62-28=34')\n    ,@('49+6=55', '20-2=18')\n    ,@('87-0=87', '82-5=77')\n    ,@('80-29=51', '51-33=18')\n    ,@('67+19=86', '42-38=4')\n    ,@('91-51=40', '15+27=42')\n    ,@('50-4=46', '5+59=64')\n    ,@('84-0=84', '13+66=79')\n    ,@('22+17=39', '36+31=67')\n    ,@('88-59=29', '67+30=97')\n    ,@('18-17=1', '23+62=85')\n    ,@('73+26=99', '86-37=49')\n    ,@('18+50=68', '58+38=96')\n    ,@('71-70=1', '34-7=27')\n    ,@('36+44=80', '37+34=71')\n    ,@('25-4=21', '38+25=63')\n    ,@('12+76=88', '24+60=84')\n    ,@('60-14=46', '83-70=13')\n    ,@('65+5=70', '4+57=61')\n    ,@('36+4=40', '78-0=78')\n    ,@('24+39=63', '26+38=64')\n    ,@('46-45=1', '90-22=68')\n    ,@('73-47=26', '89+1=90')\n    ,@('61-59=2', '88-70=18')\n    ,@('45-6=39', '3+66=69')\n    ,@('12-10=2', '85-72=13')\n    ,@('0+64=64', '9+62=71')\n    ,@('60-38=22', '77+18=95')\n    ,@('72-51=21', '7+71=78')\n    ,@('9+81=90', '10-9=1')\n    ,@('52+3=55', '72-29=43')\n    ,@('20-6=14', '22+44=66')\n    ,@('73-30=43', '53-29=24')\n    ,@('32-24=8', '46+36=82')\n    ,@('32-25=7', '24+18=42')\n    ,@('14+19=33', '41+27=68')\n    ,@('76-64=12', '90-15=75')\n    ,@('73+2=75', '20+71=91')\n    ,@('42-32=10', '25+34=59')\n    ,@('62+21=83', '3+81=84')\n    ,@('87-56=31', '49-49=0')\n    ,@('23-6=17', '74-19=55')\n    ,@('38+24=62', '51+28=79')\n    ,@('98-98=0', '34-16=18')\n    ,@('58-8=50', '99-23=76')\n    ,@('45-8=37', '98-76=22')\n    ,@('52-22=30', '90-80=10')\n    ,@('98-68=30', '40-30=10')\n    ,@('98-8=90', '97-33=64')\n    ,@('50+14=64', '13+52=65')\n    ,@('2+30=32', '68-10=58')\n    ,@('7+58=65', '75-13=62')\n    ,@('43-31=12', '40-23=17')\n    ,@('68-32=36', '22+18=40')\n    ,@('99-93=6', '51+24=75')\n    ,@('49-31=18', '59-4=55')\n    ,@('48+14=62', '33+24=57')\n    ,@('54-1=53', '18+67=85')\n    ,@('24+16=40', '72+25=97')\n    ,@('73-17=56', '21-16=5')\n    ,@('72-56=16', '71-4=67')\n    ,@('62-54=8', '46+47=93')\n    ,@('31-14=17', '12+65=77')\n    ,@('27+36=63', '95-42=53')\n    ,@('5+34=39', '30+15=45')\n    ,@('25+22=47', '12+16=28')\n    ,@('84-30=54', '97-95=2')\n    ,@('51-19=32', '85-27=58')\n    ,@('20-19=1', '83+14=97')\n    ,@('32+11=43', '14+42=56')\n    ,@('78+10=88', '3+12=15')\n    ,@('4+83=87', '83-56=27')\n    ,@('5+62=67', '96-83=13')\n    ,@('47+36=83', '68+21=89')\n    ,@('13+53=66', '23+29=52')\n    ,@('20+77=97', '51-12=39')\n    ,@('67-21=46', '39-8=31')\n    ,@('4+40=44', '95-45=50')\n    ,@('81-47=34', '43+8=51')\n    ,@('72-55=17', '8+71=79')\n    ,@('82-16=66', '29+51=80')\n    ,@('47-14=33', '10+24=34')\n    ,@('78-65=13', '94-6=88')\n    ,@('84-40=44', '99-37=62')\n    ,@('63-34=29', '13+15=28')\n    ,@('4+20=24', '79+12=91')\n    ,@('45-42=3', '80-11=69')\n    ,@('77+16=93', '21+60=81')\n    ,@('99-24=75', '25+14=39')\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$null,[ref]$true,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$true,[ref]$null,[ref]$null,[ref]$null,[ref]2) | Out-Null\n}"}
